$d = $word.ActiveDocument

# Locate the end of the run that contains "A simple demonstration of a"
# (start of the document / first paragraph) -- this is where the new
# "M2Doc version mismatch" warning block needs to be inserted, mirroring
# the "<--- <message>" warning blocks already used elsewhere in this
# template for validation errors.
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$anchor.Find.Execute("A simple demonstration of a", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertStart = $anchor.End

$arrowText = "<---"
$warningText = "M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0"

$ins = $d.Range($insertStart, $insertStart)
$ins.InsertAfter("    " + $arrowText + $warningText + "    ")

$arrowStart = $insertStart + 4
$arrowEnd = $arrowStart + $arrowText.Length
$msgStart = $arrowEnd
$msgEnd = $msgStart + $warningText.Length

$arrowRange = $d.Range($arrowStart, $arrowEnd)
$arrowRange.Font.Color = 42495
$arrowRange.Font.Size = 16
$arrowRange.HighlightColorIndex = 16

$msgRange = $d.Range($msgStart, $msgEnd)
$msgRange.Font.Color = 42495
$msgRange.Font.Size = 16
$msgRange.HighlightColorIndex = 16

Write-Output "Paragraph 1 now reads: $($d.Paragraphs(1).Range.Text)"
